$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author added a new "证券中心" (securities center) board tile right
# after row 20 ("苏州"). This pushes every following row (old 21-38,
# "厦门" .. "钓鱼岛") down by one, which is exactly what EntireRow
# insertion on row 21 gives us.
$ws.Rows.Item(21).Insert()

# Populate the freshly inserted row 21 with the new tile's data. Like the
# other "event" tiles (建筑公司, 电视台, 污水处理厂, ...) only A/B/M are
# used - no block/fee/buy/sell/pledge/upgrade numbers.
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "证券中心"
$ws.Range("M21").Value = "获得500元，然后额外获得你拥有投资项目数量*500元的奖励。"

# The unused columns between B and M shouldn't carry any leftover
# formatting/content copied down from the row above during the insert.
$ws.Range("C21:L21").ClearContents()

# Match the wrapped, vertically centered look used by the other
# long-text "public"/"project" description cells (style index 2).
$ws.Range("M21").WrapText = $true
$ws.Range("M21").VerticalAlignment = -4108

# The row grows tall enough to show the wrapped description, same as the
# other two-line description rows (28.5pt).
$ws.Rows.Item(21).RowHeight = 28.5

# Reflect where the author ended up after the edit.
$null = $ws.Range("B21").Select()
